# "Generate Report for Handback" — mark the two handed-off files as handed
# back (in sync with en-US), and fill in the Latest Target File / Latest
# Handback File / Latest Handback DateTime columns on the per-locale sheets.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$mdUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fc28ae4ee5a9d616e3fb8a4280c7f7b26e0bac82/e2e/2df024e5-f384-4970-a0a5-31d3bfdf89a3.md"
$mdName1 = "2df024e5-f384-4970-a0a5-31d3bfdf89a3.md"
$mdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fc28ae4ee5a9d616e3fb8a4280c7f7b26e0bac82/e2e/e49630ad-40ac-44e9-a236-cf67916150cf.md"
$mdName2 = "e49630ad-40ac-44e9-a236-cf67916150cf.md"

# Overview sheet: per-locale status columns (E = zh-cn, F = de-de)
$overview = $wb.Worksheets("Overview")
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText

function Set-HandbackRow {
    param($ws, $statusCell, $targetCell, $targetLinkUrl, $targetLinkName, $handbackCell, $handbackFile, $dateCell, $dateValue)

    $ws.Range($statusCell).Value = $statusText

    $ws.Range($targetCell).Value = $targetLinkName
    $ws.Hyperlinks.Add($ws.Range($targetCell), $targetLinkUrl, "", "", $targetLinkName)
    $ws.Range($targetCell).Font.Name = "Calibri"
    $ws.Range($targetCell).Font.Color = 15570276
    $ws.Range($targetCell).Font.Underline = 2

    $ws.Range($handbackCell).Value = $handbackFile

    $ws.Range($dateCell).Value = $dateValue
}

# zh-cn sheet
$zhcn = $wb.Worksheets("zh-cn")
Set-HandbackRow $zhcn "C2" "I2" $mdUrl1 $mdName1 "J2" "2df024e5-f384-4970-a0a5-31d3bfdf89a3.aff411bc4df4651b26a6202ead90ebf3565ad795.zh-cn.xlf" "K2" "2016-10-21 01:10:57"
Set-HandbackRow $zhcn "C3" "I3" $mdUrl2 $mdName2 "J3" "e49630ad-40ac-44e9-a236-cf67916150cf.cf4bac4afc76e997399249f1be0b02e3d143e5f0.zh-cn.xlf" "K3" "2016-10-21 01:10:57"

# de-de sheet
$dede = $wb.Worksheets("de-de")
Set-HandbackRow $dede "C2" "I2" $mdUrl1 $mdName1 "J2" "2df024e5-f384-4970-a0a5-31d3bfdf89a3.aff411bc4df4651b26a6202ead90ebf3565ad795.de-de.xlf" "K2" "2016-10-21 01:11:16"
Set-HandbackRow $dede "C3" "I3" $mdUrl2 $mdName2 "J3" "e49630ad-40ac-44e9-a236-cf67916150cf.cf4bac4afc76e997399249f1be0b02e3d143e5f0.de-de.xlf" "K3" "2016-10-21 01:11:16"

Write-Host "Handback report generated."
